# "GET TESTS - Share test"
# Re-color existing checklist rows to reflect updated task status, clear the
# now-obsolete "Matilda" owner notes on a couple of rows, assign "Philip" as
# owner on a few items, and append a new block of checklist rows (a divider
# plus new GET TESTS / share-test follow-up tasks), some owned by Markus.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$AlignCenter    = -4108
$AlignBottom    = -4107

# ---- Color constants (OLE BGR-packed RGB, i.e. R + G*256 + B*65536) ----
$Green00B050 = 5287936  # FF00B050 (new "done" color)
$Gray7F7F7F  = 8355711  # new neutral divider color
$GreenTheme  = 4697456  # FF70AD47 (same visual color as the sheet's theme9 accent)
$Red         = 255      # FFFF0000

# ---- Re-style existing rows (status colors changed), by cloning the
#      formatting of an already-present cell with the matching look, then
#      nudging fill/alignment where the target differs slightly. This keeps
#      reusing the workbook's existing theme-based fills instead of minting
#      new literal-RGB duplicates. ----

$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial($xlPasteFormats)
$ws.Range("A4").VerticalAlignment = $AlignBottom

$ws.Range("A9").Copy()
$ws.Range("A5").PasteSpecial($xlPasteFormats)

$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial($xlPasteFormats)

$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial($xlPasteFormats)
$ws.Range("A7").VerticalAlignment = $AlignBottom

$ws.Range("A9").Copy()
$ws.Range("A8").PasteSpecial($xlPasteFormats)

$ws.Range("A9").PasteSpecial($xlPasteFormats)
$ws.Range("A9").Interior.Color = $Green00B050

$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial($xlPasteFormats)

$ws.Range("A5").Copy()
$ws.Range("A14").PasteSpecial($xlPasteFormats)

$ws.Range("A9").Copy()
$ws.Range("A15").PasteSpecial($xlPasteFormats)

# ---- Owner column updates ----
$ws.Range("B7").ClearContents()
$ws.Range("B9").ClearContents()

$ws.Range("B8").Value = "Philip"
$ws.Range("B10").Value = "Philip"
$ws.Range("B15").Value = "Philip"

# ---- New rows: blank divider + new checklist items ----

# Row 16: empty divider cell - Times New Roman, centered, gray fill.
$ws.Range("A13").Copy()
$ws.Range("A16").PasteSpecial($xlPasteFormats)
$ws.Range("A16").Interior.Color = $Gray7F7F7F

# Rows 17/18: default font, theme-green fill, no explicit alignment.
# Build the combo once in a scratch cell (off the used range) and clone it,
# so both rows share a single new fill entry instead of minting two.
$ws.Range("Z1").Interior.Color = $GreenTheme
$ws.Range("Z1").Copy()

$ws.Range("A17").Value = "Ändra alla knappar, labels,textfield. De ska skapas av klassen CreateNodes"
$ws.Range("A17").PasteSpecial($xlPasteFormats)

$ws.Range("A18").Value = "Shuffla svarsalternativ när användaren skriver ett prov"
$ws.Range("A18").PasteSpecial($xlPasteFormats)
$ws.Range("B18").Value = "Markus"

$ws.Range("Z1").Clear()

# Rows 19/20: default font, red fill, no explicit alignment.
$ws.Range("Z1").Interior.Color = $Red
$ws.Range("Z1").Copy()

$ws.Range("A19").Value = "Ta bort navbar vid utloggning"
$ws.Range("A19").PasteSpecial($xlPasteFormats)

$ws.Range("A20").Value = "kommentera kod"
$ws.Range("A20").PasteSpecial($xlPasteFormats)

$ws.Range("Z1").Clear()

# ---- Selection as left by the author after the edit ----
$ws.Range("A23").Select() | Out-Null
